# Generate Report for Handoff
# Replaces the stale GUID-named markdown file (3984c643-...) with a freshly
# generated handoff (4648d857-...), refreshes the handoff timestamps/xlf
# names, and clears the (not-yet-populated) handback columns now that the
# file is freshly re-handed-off.

$wb = $excel.ActiveWorkbook

$newGuidFile   = "4648d857-fa1a-4d25-a2c5-ac4202fb41ea.md"
$newGuidPath   = "e2e\" + $newGuidFile

$newXliffBase  = "4648d857-fa1a-4d25-a2c5-ac4202fb41ea.a6dac6e71af063ebd498eb5223322dccdb5fb23e"
$newXliffZhCn  = $newXliffBase + ".zh-cn.xlf"
$newXliffDeDe  = $newXliffBase + ".de-de.xlf"

$newOverviewDate = "2016-08-12 19:14:50"
$newZhCnHoDate    = "2016-08-12 19:14:42"
$emptyHandback    = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")

$ovw.Range("A2").Value = $newGuidFile
$ovw.Range("B2").Value = $newGuidPath
$ovw.Range("G2").Value = $newOverviewDate

foreach ($hl in $ovw.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.TextToDisplay = $newGuidPath
    }
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = $newGuidFile
$zh.Range("G2").Value = $newXliffZhCn
$zh.Range("H2").Value = $newZhCnHoDate

# No handback has happened yet for the freshly generated handoff -
# clear the stale target/handback file columns and reset the handback
# datetime to the "never" sentinel.
$zh.Range("I2").Value = ""
$zh.Range("J2").Value = ""
$zh.Range("K2").Value = $emptyHandback

foreach ($hl in $zh.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $newGuidFile
    } elseif ($addr -eq '$I$2') {
        $hl.Delete()
    }
}

$zh.Range("I2").Style = "Normal"

$zh.Columns.Item(9).ColumnWidth = 17.8
$zh.Columns.Item(10).ColumnWidth = 20.8

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = $newGuidFile
$de.Range("G2").Value = $newXliffDeDe
$de.Range("H2").Value = $newOverviewDate
$de.Range("I2").Value = ""
$de.Range("J2").Value = ""
$de.Range("K2").Value = $emptyHandback

foreach ($hl in $de.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $newGuidFile
    } elseif ($addr -eq '$I$2') {
        $hl.Delete()
    }
}

$de.Range("I2").Style = "Normal"

$de.Columns.Item(9).ColumnWidth = 17.8
$de.Columns.Item(10).ColumnWidth = 20.8
